$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = "work in progress"
$ws.Range("F6").Value = "Ended"
$ws.Range("G6").Value = "home made encoder / at the end work with analog read and asynchronous soft"
$ws.Range("I6").Value = "not that much precise regarding rotation"
$ws.Range("H9").Value = "Tested:  arduino connection and dialog >> decided to just use magnetometer as a first step"
$ws.Range("E12").Value = "december 2015"
$ws.Range("E14").Value = "january 2016"
$ws.Range("H14").Value = "ok"
$ws.Range("I14").Value = "`nto be design as an object "
$ws.Range("I15").Value = "encoders & echo system > permanent power consumption "
$ws.Range("B16").Value = "Integration Octave"
$ws.Range("C16").Value = "determine the possibilities to integrate octave and java code"
$ws.Range("H16").Value = "java defined as octave object "
$ws.Range("B17").Value = "Asynchronous communication"
$ws.Range("C17").Value = "developp an asynchronous communication between Octave and Arduino using Java server"
$ws.Range("C18").Value = "add northorientation to the learning data of the scans"
$ws.Range("F18").Value = "canceled"
$ws.Range("A19").Value = "Development"
$ws.Range("B19").Value = "Echo Localization"
$ws.Range("C19").Value = "developp a flat learning and use North Orientation to shift current scan according to flat measurment"
$ws.Range("F19").Value = "Ended"
$ws.Range("A21").Value = "Development"
$ws.Range("B21").Value = "Define a global loop logic"
$ws.Range("C21").Value = "Octave code"
$ws.Range("F21").Value = "Started"
$ws.Range("A22").Value = "Integration"
$ws.Range("B22").Value = "Design an eletronic power switch "
$ws.Range("C22").Value = "to reduce power consumption of encoders & echo system"
$ws.Range("A24").Value = "Development"
$ws.Range("B24").Value = "Check constitency between robot heading calculation and North Orientation"
$ws.Range("F24").Value = "Started"
$ws.Range("G24").Value = "octave loop"
$ws.Range("A25").Value = "Development"
$ws.Range("B25").Value = "Path determination"
$ws.Range("C25").Value = "octave Astar code"
$ws.Range("F25").Value = "Started"
$ws.Range("A26").Value = "Development"
$ws.Range("B26").Value = "Use mode pulse & NO to rotate little angle"
$ws.Range("C26").Value = "arduino code"
$ws.Range("A27").Value = "Development"
$ws.Range("B27").Value = "Use NO to imporve rotation, precision"
$ws.Range("C27").Value = "arduino code"
$ws.Range("A28").Value = "Development"
$ws.Range("B28").Value = "Use echo F B before and after move"
$ws.Range("C28").Value = "to increase location precision"
$ws.Range("A29").Value = "Development"
$ws.Range("B29").Value = "Astar Path detemination"
$ws.Range("C29").Value = "add possibility to inhibit some kind of move according to obastclea"
$ws.Range("A30").Value = "Development"
$ws.Range("B30").Value = "Use echo F B Right Left before rotation"
$ws.Range("C30").Value = "to be sure to avoid osbacle"
$ws.Range("A31").Value = "Development"
$ws.Range("B31").Value = "Path Calculation add logic to take into account robot shape"
$ws.Range("A32").Value = "Development"
$ws.Range("B32").Value = "Octave developp a localization convergence loop"
$ws.Range("C32").Value = "include robot move"

# Column B width
$ws.Columns.Item(2).ColumnWidth = 57

# Selection matches final cursor position
$ws.Range("G34").Select()
